$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$d = $word.ActiveDocument

# 1. Paragraph 1 ("Rekayasa Perangkat Lunak"): drop the _GoBack bookmark that
#    used to sit here; it moves down to the new "Smester" paragraph below.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML('<w:p ' + $ns + '><w:r><w:t>Rekayasa Perangkat Lunak</w:t></w:r></w:p>')

# 2. Paragraph 2 ("Nama ... : Dindin Sihabudin Ahmad Wardi"): the single
#    "<w:tab/><w:t>...</w:t>" run becomes two tab runs followed by the text run.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML('<w:p ' + $ns + '><w:r><w:t xml:space="preserve">Nama </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>: Dindin Sihabudin Ahmad Wardi</w:t></w:r></w:p>')

# 3. Paragraph 3 ("NPM ... : 1406037"): same tab-run split as paragraph 2.
$p3 = $d.Paragraphs(3)
$p3.Range.InsertXML('<w:p ' + $ns + '><w:r><w:t>NPM</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>: 1406037</w:t></w:r></w:p>')

# 4. New paragraph 4 ("Smester ... : VI (Enam)") inserted after paragraph 3,
#    carrying the _GoBack bookmark that was removed from paragraph 1.
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.InsertXML('<w:p ' + $ns + '><w:r><w:t>Smester</w:t></w:r><w:r><w:tab/><w:t>: VI (Enam)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
